$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 5235.4287
$ws.Range("I62").Value = 4869.6
$ws.Range("K62").Value = 4869.6
$ws.Range("M62").Value = -4245.6
$ws.Range("H65").Value = 5235.4287
$ws.Range("I65").Value = 4869.6
$ws.Range("K65").Value = 24348
$ws.Range("M65").Value = -21228
$ws.Range("H80").Value = 2472.3513
$ws.Range("I80").Value = 2512.3157
$ws.Range("J80").Value = 2430.1667
$ws.Range("K80").Value = 7536.9471
$ws.Range("L80").Value = 7290.500100000001
$ws.Range("M80").Value = -6538.9471
$ws.Range("N80").Value = -9286.500100000001
$ws.Range("H83").Value = 2472.3513
$ws.Range("I83").Value = 2512.3157
$ws.Range("J83").Value = 2430.1667
$ws.Range("K83").Value = 22610.8413
$ws.Range("L83").Value = 21871.5003
$ws.Range("M83").Value = -17618.8413
$ws.Range("N83").Value = -31855.5003
$ws.Range("H86").Value = 2078.1724
$ws.Range("J86").Value = 2843.2856
$ws.Range("L86").Value = 2843.2856
$ws.Range("N86").Value = -5089.2856
$ws.Range("H88").Value = 2258.6667
$ws.Range("I88").Value = 2510.4
$ws.Range("J88").Value = 1000
$ws.Range("K88").Value = 2510.4
$ws.Range("L88").Value = 1000
$ws.Range("M88").Value = -2104.4
$ws.Range("N88").Value = -1812
$ws.Range("H89").Value = 2078.1724
$ws.Range("J89").Value = 2843.2856
$ws.Range("L89").Value = 14216.428
$ws.Range("N89").Value = -25448.428
$ws.Range("H91").Value = 2258.6667
$ws.Range("I91").Value = 2510.4
$ws.Range("J91").Value = 1000
$ws.Range("K91").Value = 2510.4
$ws.Range("L91").Value = 1000
$ws.Range("M91").Value = -1106.4
$ws.Range("N91").Value = -3808
$ws.Range("H135").Value = 13159679
$ws.Range("I135").Value = 1401.091
$ws.Range("K135").Value = 12609.819
$ws.Range("M135").Value = -10074.819
$ws.Range("H138").Value = 7410624.5
$ws.Range("J138").Value = 9012385
$ws.Range("L138").Value = 27037155
$ws.Range("N138").Value = -27047435
$ws.Range("H141").Value = 2775.5334
$ws.Range("I141").Value = 2823.4285
$ws.Range("J141").Value = 2105
$ws.Range("K141").Value = 8470.2855
$ws.Range("L141").Value = 6315
$ws.Range("M141").Value = -3290.2855
$ws.Range("N141").Value = -16675

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 14495376
$ws.Range("I61").Value = 17243366
$ws.Range("K61").Value = 17243366
$ws.Range("M61").Value = -17243154
$ws.Range("H74").Value = 47675980
$ws.Range("I74").Value = 66743560
$ws.Range("J74").Value = 7031.1665
$ws.Range("K74").Value = 66743560
$ws.Range("L74").Value = 7031.1665
$ws.Range("M74").Value = -66742686
$ws.Range("N74").Value = -8779.166499999999
$ws.Range("H77").Value = 47675980
$ws.Range("I77").Value = 66743560
$ws.Range("J77").Value = 7031.1665
$ws.Range("K77").Value = 333717800
$ws.Range("L77").Value = 35155.8325
$ws.Range("M77").Value = -333713432
$ws.Range("N77").Value = -43891.8325
$ws.Range("H136").Value = 14495376
$ws.Range("I136").Value = 17243366
$ws.Range("K136").Value = 51730098
$ws.Range("M136").Value = -51727548

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1423.125
$ws.Range("I94").Value = 1423.125
$ws.Range("K94").Value = 1423.125
$ws.Range("M94").Value = -972.125

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H74").Value = 43095
$ws.Range("H77").Value = 43095
$ws.Range("H99").Value = 7244.5835
$ws.Range("J99").Value = 3839.4
$ws.Range("L99").Value = 3839.4
$ws.Range("N99").Value = -6835.4
$ws.Range("H106").Value = 63000
$ws.Range("J106").Value = 63000
$ws.Range("L106").Value = 63000
$ws.Range("N106").Value = -65524
$ws.Range("H126").Value = 7244.5835
$ws.Range("J126").Value = 3839.4
$ws.Range("L126").Value = 11518.2
$ws.Range("N126").Value = -16458.2
$ws.Range("H131").Value = 85733
$ws.Range("J131").Value = 85733
$ws.Range("L131").Value = 85733
$ws.Range("N131").Value = -95813
$ws.Range("H141").Value = 495696.9
$ws.Range("J141").Value = 590407.3
$ws.Range("L141").Value = 590407.3
$ws.Range("N141").Value = -600767.3

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 2050626.8
$ws.Range("I2").Value = 2344.8572
$ws.Range("J2").Value = 4440289
$ws.Range("K2").Value = 14069.1432
$ws.Range("L2").Value = 26641734
$ws.Range("M2").Value = -13956.1432
$ws.Range("N2").Value = -26641960
$ws.Range("H6").Value = 51.92857
$ws.Range("I6").Value = 51.92857
$ws.Range("K6").Value = 155.78571
$ws.Range("M6").Value = -42.78570999999999
$ws.Range("H34").Value = 1860.3
$ws.Range("I34").Value = 287.9524
$ws.Range("J34").Value = 3598.158
$ws.Range("K34").Value = 863.8572
$ws.Range("L34").Value = 10794.474
$ws.Range("M34").Value = -779.8572
$ws.Range("N34").Value = -10962.474
$ws.Range("H107").Value = 862.44446
$ws.Range("I107").Value = 531.6
$ws.Range("J107").Value = 1276
$ws.Range("K107").Value = 1594.8
$ws.Range("L107").Value = 3828
$ws.Range("M107").Value = 325.1999999999998
$ws.Range("N107").Value = -7668
$ws.Range("H128").Value = 121575
$ws.Range("I128").Value = 121575
$ws.Range("K128").Value = 364725
$ws.Range("M128").Value = -359745

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 7902794
$ws.Range("I14").Value = 11286849
$ws.Range("J14").Value = 6666.6665
$ws.Range("K14").Value = 11286849
$ws.Range("L14").Value = 6666.6665
$ws.Range("M14").Value = -11286681
$ws.Range("N14").Value = -7002.6665
$ws.Range("H126").Value = 14291696
$ws.Range("I126").Value = 10007818
$ws.Range("K126").Value = 30023454
$ws.Range("M126").Value = -30020984

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4444
$ws.Range("I40").Value = 4143.8823
$ws.Range("K40").Value = 4143.8823
$ws.Range("M40").Value = -4007.8823
$ws.Range("H93").Value = 3500
$ws.Range("J93").Value = 3500
$ws.Range("L93").Value = 3500
$ws.Range("N93").Value = -5996
$ws.Range("H122").Value = 5134.9
$ws.Range("I122").Value = 4594.3335
$ws.Range("J122").Value = 10000
$ws.Range("K122").Value = 13783.0005
$ws.Range("L122").Value = 30000
$ws.Range("M122").Value = -11333.0005
$ws.Range("N122").Value = -34900
$ws.Range("H136").Value = 4387.2354
$ws.Range("I136").Value = 4387.2354
$ws.Range("K136").Value = 13161.7062
$ws.Range("M136").Value = -10611.7062
